$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "2025/12/06 02:00"
$ws.Range("B80").Value = "-"
$ws.Range("C80").Value = "-"
$ws.Range("D80").Value = "-"
$ws.Range("E80").Value = "-"
$ws.Range("F80").Value = "-"
$ws.Range("G80").Value = "-"
